$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (rich-text cells, set as flattened text) ---
$ws.Range("A8").Value = "Volume 32   Number  4"
$ws.Range("C9").Value = "Report Covering the Week  1/20/2025  Through  1/26/2025"

# --- Reference/template cells used to restore the correct style index after a
#     number<->text type change (Value= alone resets styles on type flips). ---
# Style 13 = text cell (General fmt, right aligned)
# Style 14 = numeric, "#,##0.0;-#,##0.0" (pct-change columns)
# Style 15 = numeric, "#,##0" (count columns)
$refText = $ws.Range("D15")   # style 13 the whole edit through
$refPct  = $ws.Range("H15")   # style 14 the whole edit through
$refNum  = $ws.Range("F15")   # style 15 the whole edit through

function Set-Cell($ws, $addr, $value, $restyleRef) {
    $cell = $ws.Range($addr)
    if ($restyleRef -ne $null -and $value -is [string]) {
        # Force text storage, then restore the canonical style via PasteSpecial(formats).
        $cell.NumberFormat = "@"
        $cell.Value = $value
        $restyleRef.Copy()
        $cell.PasteSpecial(-4122)
    } elseif ($restyleRef -ne $null) {
        $cell.Value = $value
        $restyleRef.Copy()
        $cell.PasteSpecial(-4122)
    } else {
        $cell.Value = $value
    }
}

# Row 14
Set-Cell $ws "N14" -100 $refPct

# Row 15
Set-Cell $ws "C15" "0" $refText

# Row 16
Set-Cell $ws "D16" 6 $null
Set-Cell $ws "E16" -66.666666666666 $null
Set-Cell $ws "F16" 8 $null
Set-Cell $ws "G16" 19 $null
Set-Cell $ws "H16" -57.894736842105 $null
Set-Cell $ws "I16" 7 $null
Set-Cell $ws "J16" 16 $null
Set-Cell $ws "K16" -56.25 $null
Set-Cell $ws "L16" -36.363636363636 $null
Set-Cell $ws "M16" -22.222222222222 $null
Set-Cell $ws "N16" -86.538461538461 $null

# Row 17
Set-Cell $ws "C17" 6 $null
Set-Cell $ws "D17" 3 $null
Set-Cell $ws "E17" 100 $null
Set-Cell $ws "F17" 17 $null
Set-Cell $ws "G17" 20 $null
Set-Cell $ws "H17" -15 $null
Set-Cell $ws "I17" 17 $null
Set-Cell $ws "J17" 17 $null
Set-Cell $ws "K17" 0 $null
Set-Cell $ws "L17" 54.545454545454 $null
Set-Cell $ws "M17" 142.857142857143 $null
Set-Cell $ws "N17" 54.545454545454 $null

# Row 18
Set-Cell $ws "D18" 1 $null
Set-Cell $ws "E18" 100 $null
Set-Cell $ws "I18" 8 $null
Set-Cell $ws "J18" 7 $null
Set-Cell $ws "K18" 14.285714285714 $null
Set-Cell $ws "L18" -11.111111111111 $null
Set-Cell $ws "M18" -57.894736842105 $null
Set-Cell $ws "N18" -90.47619047619 $null

# Row 19
Set-Cell $ws "C19" 2 $null
Set-Cell $ws "D19" 9 $null
Set-Cell $ws "E19" -77.777777777777 $null
Set-Cell $ws "F19" 32 $null
Set-Cell $ws "G19" 54 $null
Set-Cell $ws "H19" -40.74074074074 $null
Set-Cell $ws "I19" 28 $null
Set-Cell $ws "J19" 53 $null
Set-Cell $ws "K19" -47.169811320754 $null
Set-Cell $ws "L19" -31.70731707317 $null
Set-Cell $ws "M19" 40 $null
Set-Cell $ws "N19" 27.272727272727 $null

# Row 20
Set-Cell $ws "C20" 6 $null
Set-Cell $ws "D20" 8 $null
Set-Cell $ws "E20" -25 $null
Set-Cell $ws "F20" 14 $null
Set-Cell $ws "G20" 36 $null
Set-Cell $ws "H20" -61.111111111111 $null
Set-Cell $ws "I20" 14 $null
Set-Cell $ws "J20" 33 $null
Set-Cell $ws "K20" -57.575757575757 $null
Set-Cell $ws "L20" -48.148148148148 $null
Set-Cell $ws "M20" 100 $null
Set-Cell $ws "N20" -88.709677419354 $null

# Row 21
Set-Cell $ws "C21" 18 $null
Set-Cell $ws "D21" 27 $null
Set-Cell $ws "E21" -33.333333333333 $null
Set-Cell $ws "F21" 83 $null
Set-Cell $ws "G21" 138 $null
Set-Cell $ws "H21" -39.855072463768 $null
Set-Cell $ws "I21" 76 $null
Set-Cell $ws "J21" 127 $null
Set-Cell $ws "K21" -40.15748031496 $null
Set-Cell $ws "L21" -23.232323232323 $null
Set-Cell $ws "M21" 20.63492063492 $null
Set-Cell $ws "N21" -74.581939799331 $null

# Row 22
Set-Cell $ws "L22" -100 $refPct

# Row 23
Set-Cell $ws "C23" "0" $refText
Set-Cell $ws "D23" 4 $refNum
Set-Cell $ws "E23" -100 $refPct
Set-Cell $ws "F23" 5 $null
Set-Cell $ws "G23" 7 $null
Set-Cell $ws "H23" -28.571428571428 $null
Set-Cell $ws "J23" 7 $null
Set-Cell $ws "K23" -28.571428571428 $null
Set-Cell $ws "M23" 25 $null

# Row 24
Set-Cell $ws "C24" 22 $null
Set-Cell $ws "D24" 15 $null
Set-Cell $ws "E24" 46.666666666666 $null
Set-Cell $ws "F24" 66 $null
Set-Cell $ws "G24" 78 $null
Set-Cell $ws "H24" -15.384615384615 $null
Set-Cell $ws "I24" 62 $null
Set-Cell $ws "J24" 71 $null
Set-Cell $ws "K24" -12.676056338028 $null
Set-Cell $ws "L24" -18.421052631578 $null
Set-Cell $ws "M24" 47.619047619047 $null

# Row 25
Set-Cell $ws "C25" 3 $null
Set-Cell $ws "D25" 6 $null
Set-Cell $ws "E25" -50 $null
Set-Cell $ws "F25" 16 $null
Set-Cell $ws "G25" 21 $null
Set-Cell $ws "H25" -23.809523809523 $null
Set-Cell $ws "I25" 15 $null
Set-Cell $ws "J25" 21 $null
Set-Cell $ws "K25" -28.571428571428 $null
Set-Cell $ws "L25" -61.538461538461 $null

# Row 26
Set-Cell $ws "D26" 4 $null
Set-Cell $ws "E26" 125 $null
Set-Cell $ws "F26" 28 $null
Set-Cell $ws "G26" 19 $null
Set-Cell $ws "H26" 47.368421052631 $null
Set-Cell $ws "I26" 26 $null
Set-Cell $ws "J26" 17 $null
Set-Cell $ws "K26" 52.941176470588 $null
Set-Cell $ws "L26" 13.043478260869 $null
Set-Cell $ws "M26" -7.142857142857 $null

# Row 27
Set-Cell $ws "C27" "0" $refText
Set-Cell $ws "E27" -100 $null
Set-Cell $ws "G27" 4 $null
Set-Cell $ws "H27" -25 $null
Set-Cell $ws "J27" 4 $null
Set-Cell $ws "K27" -25 $null

# Row 28
Set-Cell $ws "C28" 1 $null

# Row 31
Set-Cell $ws "F31" "0" $refText
Set-Cell $ws "G31" "0" $refText
Set-Cell $ws "H31" "***.*" $refText
